# Auto-generated edit script applying the crypto price/volume refresh
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new value is a plain decimal number (e.g. "247.48") must be
# forced to Text format first, otherwise Excel auto-converts them to numeric
# cells -- the source data stores every Price/Volume cell as text.
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D51').NumberFormat = '@'

$ws.Range('D2').Value = '91.965.26'
$ws.Range('E2').Value = '  +0.94%  '
$ws.Range('D3').Value = '3.120.98'
$ws.Range('E3').Value = '  -0.68%  '
$ws.Range('E4').Value = '  +0.06%  '
$ws.Range('D5').Value = '247.48'
$ws.Range('E5').Value = '  +2.95%  '
$ws.Range('D6').Value = '623.78'
$ws.Range('E6').Value = '  -2.00%  '
$ws.Range('E7').Value = '  +6.19%  '
$ws.Range('D8').Value = '0.373'
$ws.Range('E8').Value = '  +1.94%  '
$ws.Range('E9').Value = '  -0.04%  '
$ws.Range('D10').Value = '3.117.66'
$ws.Range('E10').Value = '  -0.70%  '
$ws.Range('D11').Value = '0.761'
$ws.Range('E11').Value = '  +5.35%  '
$ws.Range('E12').Value = '  +3.36%  '
$ws.Range('E13').Value = '  +2.35%  '
$ws.Range('D14').Value = '35.70'
$ws.Range('E14').Value = '  -3.56%  '
$ws.Range('D15').Value = '91.673.69'
$ws.Range('E15').Value = '  +0.84%  '
$ws.Range('D16').Value = '5.50'
$ws.Range('E16').Value = '  -0.92%  '
$ws.Range('D17').Value = '3.695.81'
$ws.Range('E17').Value = '  -0.68%  '
$ws.Range('D18').Value = '3.143.84'
$ws.Range('E18').Value = '  +0.39%  '
$ws.Range('D19').Value = '3.76'
$ws.Range('E19').Value = '  -0.32%  '
$ws.Range('D20').Value = '14.68'
$ws.Range('E20').Value = '  +2.62%  '
$ws.Range('D21').Value = '0.0000216'
$ws.Range('E21').Value = '  +1.61%  '
$ws.Range('D22').Value = '5.82'
$ws.Range('E22').Value = '  +3.10%  '
$ws.Range('D23').Value = '448.98'
$ws.Range('E23').Value = '  +0.61%  '
$ws.Range('D24').Value = '9.34'
$ws.Range('E24').Value = '  +3.32%  '
$ws.Range('D25').Value = '5.93'
$ws.Range('E25').Value = '  -0.78%  '
$ws.Range('D26').Value = '91.19'
$ws.Range('E26').Value = '  +1.01%  '
$ws.Range('D27').Value = '12.10'
$ws.Range('E27').Value = '  -3.18%  '
$ws.Range('D28').Value = '3.263.22'
$ws.Range('E28').Value = '  -0.99%  '
$ws.Range('E29').Value = '  +0.02%  '
$ws.Range('D30').Value = '0.187'
$ws.Range('E30').Value = '  +16.60%  '
$ws.Range('E31').Value = '  +21.14%  '
$ws.Range('D32').Value = '9.41'
$ws.Range('E32').Value = '  -2.83%  '
$ws.Range('D33').Value = '0.174'
$ws.Range('E33').Value = '  +16.43%  '
$ws.Range('D34').Value = '1.01'
$ws.Range('E34').Value = '  +3.73%  '
$ws.Range('D35').Value = '0.111'
$ws.Range('E35').Value = '  +29.89%  '
$ws.Range('B36').Value = 'RenderToken'
$ws.Range('C36').Value = 'https://coinranking.com/coin/vfo5XYwcV+rendertoken-render'
$ws.Range('D36').Value = '7.98'
$ws.Range('E36').Value = '  +10.88%  '
$ws.Range('B37').Value = 'EthereumClassic'
$ws.Range('C37').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D37').Value = '26.71'
$ws.Range('E37').Value = '  -1.58%  '
$ws.Range('D38').Value = '4.14'
$ws.Range('E38').Value = '  +24.18%  '
$ws.Range('D39').Value = '1.93'
$ws.Range('E39').Value = '  -0.65%  '
$ws.Range('D40').Value = '496.54'
$ws.Range('E40').Value = '  -3.73%  '
$ws.Range('D41').Value = '3.66'
$ws.Range('E41').Value = '  -4.21%  '
$ws.Range('D42').Value = '1.31'
$ws.Range('E42').Value = '  +0.29%  '
$ws.Range('D43').Value = '0.426'
$ws.Range('E43').Value = '  +1.63%  '
$ws.Range('D44').Value = '22.18'
$ws.Range('E44').Value = '  -0.15%  '
$ws.Range('E45').Value = '  -0.02%  '
$ws.Range('D46').Value = '1.93'
$ws.Range('E46').Value = '  -0.89%  '
$ws.Range('D47').Value = '0.702'
$ws.Range('E47').Value = '  +1.03%  '
$ws.Range('D48').Value = '154.14'
$ws.Range('E48').Value = '  +2.31%  '
$ws.Range('D49').Value = '4.56'
$ws.Range('E49').Value = '  +0.49%  '
$ws.Range('D50').Value = '1.36'
$ws.Range('E50').Value = '  -0.89%  '
$ws.Range('D51').Value = '44.37'
$ws.Range('E51').Value = '  -3.30%  '
